# Update countries & provincias Spain
# The "Pais" sheet is a COVID-19 ranking table (rows sorted by column B,
# "Casos totales", descending). Refreshing the source data nudges a few
# countries past their neighbours in the ranking, so for those rows the
# country name in column A has to move along with the new numbers while
# the numeric snapshot (B:H) for every updated row is refreshed to the
# newer figures. The footer timestamp in A1 is bumped too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" footer text (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 26 de Agosto de 2020 a las 11:40"

# --- Countries that changed rank order (names swap between these rows) ---
$ws.Range("A47").Value = "Polonia"
$ws.Range("A48").Value = "Japon"

$ws.Range("A121").Value = "Eslovaquia"
$ws.Range("A122").Value = "Mozambique"

$ws.Range("A129").Value = "Eslovenia"
$ws.Range("A130").Value = "Mali"
$ws.Range("A131").Value = "Gambia"

$ws.Range("A169").Value = "Birmania"
$ws.Range("A170").Value = "Tanzania"

# --- Refreshed case numbers (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Row 6 (India)
$ws.Range("B6").Value = 3239096
$ws.Range("C6").Value = 7342
$ws.Range("E6").Value = 710763
$ws.Range("G6").Value = 33
$ws.Range("H6").Value = 59645

# Row 23 (Alemania)
$ws.Range("B23").Value = 237589
$ws.Range("C23").Value = 17
$ws.Range("E23").Value = 18644

# Row 26 (Indonesia)
$ws.Range("B26").Value = 160165
$ws.Range("C26").Value = 2306
$ws.Range("D26").Value = 115409
$ws.Range("E26").Value = 37812
$ws.Range("G26").Value = 86
$ws.Range("H26").Value = 6944

# Row 32 (Israel)
$ws.Range("B32").Value = 107341
$ws.Range("C32").Value = 881
$ws.Range("D32").Value = 85893
$ws.Range("E32").Value = 20581
$ws.Range("G32").Value = 8
$ws.Range("H32").Value = 867

# Row 47 (Polonia, after rename)
$ws.Range("B47").Value = 63802
$ws.Range("C47").Value = 729
$ws.Range("D47").Value = 43399
$ws.Range("E47").Value = 18409
$ws.Range("G47").Value = 17
$ws.Range("H47").Value = 1994

# Row 48 (Japon, after rename)
$ws.Range("B48").Value = 63121
$ws.Range("D48").Value = 50431
$ws.Range("E48").Value = 11494
$ws.Range("H48").Value = 1196

# Row 91 (Consejo Danes para los Refugiados)
$ws.Range("B91").Value = 9912
$ws.Range("C91").Value = 21
$ws.Range("D91").Value = 8987
$ws.Range("E91").Value = 671
$ws.Range("G91").Value = 3
$ws.Range("H91").Value = 254

# Row 92 (Malasia)
$ws.Range("B92").Value = 9291
$ws.Range("C92").Value = 6
$ws.Range("D92").Value = 8978
$ws.Range("E92").Value = 188

# Row 101 (Finlandia)
$ws.Range("B101").Value = 8002
$ws.Range("C101").Value = 21
$ws.Range("E101").Value = 567

# Row 121 (Eslovaquia, after rename)
$ws.Range("B121").Value = 3536
$ws.Range("C121").Value = 84
$ws.Range("D121").Value = 2192
$ws.Range("E121").Value = 1311
$ws.Range("H121").Value = 33

# Row 122 (Mozambique, after rename)
$ws.Range("B122").Value = 3508
$ws.Range("D122").Value = 1809
$ws.Range("E122").Value = 1678
$ws.Range("H122").Value = 21

# Row 127 (Sri Lanka)
$ws.Range("D127").Value = 2819
$ws.Range("E127").Value = 140

# Row 129 (Eslovenia, after rename)
$ws.Range("B129").Value = 2722
$ws.Range("C129").Value = 36
$ws.Range("D129").Value = 2170
$ws.Range("E129").Value = 419
$ws.Range("H129").Value = 133

# Row 130 (Mali, after rename)
$ws.Range("B130").Value = 2713
$ws.Range("D130").Value = 2041
$ws.Range("E130").Value = 547
$ws.Range("H130").Value = 125

# Row 131 (Gambia, after rename)
$ws.Range("D131").Value = 601
$ws.Range("E131").Value = 1995
$ws.Range("H131").Value = 90

# Row 169 (Birmania, after rename)
$ws.Range("B169").Value = 557
$ws.Range("C169").Value = 53
$ws.Range("D169").Value = 341
$ws.Range("E169").Value = 210
$ws.Range("H169").Value = 6

# Row 170 (Tanzania, after rename)
$ws.Range("B170").Value = 509
$ws.Range("D170").Value = 183
$ws.Range("E170").Value = 305
$ws.Range("H170").Value = 21
